$d = $word.ActiveDocument

# Locate the "Device:" table (6 columns: Device:, Laryngoscope, lma_details,
# glide_details, wis_hipple_details, other_device_details).
$t = $null
for ($i = 1; $i -le $d.Tables.Count; $i++) {
    $candidate = $d.Tables.Item($i)
    if ($candidate.Columns.Count -eq 6 -and $candidate.Cell(1, 1).Range.Text.StartsWith("Device:")) {
        $t = $candidate
    }
}

# Remove the "wis_hipple_details" column (originally the 5th column / 2119 dxa).
$t.Columns.Item(5).Delete()

# Resize the remaining data columns (2-5) to match the new layout.
$t.Cell(1, 2).Width = 2453 / 20.0
$t.Cell(1, 3).Width = 2453 / 20.0
$t.Cell(1, 4).Width = 2453 / 20.0
$t.Cell(1, 5).Width = 2454 / 20.0

# Center-align the paragraph text in each of those data cells.
$t.Cell(1, 2).Range.Paragraphs.Item(1).Alignment = 1
$t.Cell(1, 3).Range.Paragraphs.Item(1).Alignment = 1
$t.Cell(1, 4).Range.Paragraphs.Item(1).Alignment = 1
$t.Cell(1, 5).Range.Paragraphs.Item(1).Alignment = 1
